$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 ---
$ws.Range("A11").Value = 18
$ws.Range("B11").Value = 61
$ws.Range("C11").Value = "09.09.2023, 22:59:56"
$ws.Range("D11").Value = "Блок Польский"
$ws.Range("E11").Value = "туман"
$ws.Range("F11").Value = "УЦЕНКА"
$ws.Range("G11").Value = 680
$ws.Range("H11").Value = "Продукция Плиточка/Блоки заборные"

# I11 looks like a date ("08.09.2023") but must stay stored as plain text,
# so force the cell to Text format before assigning, then restore the
# default "Normal" style so no extra number-format style gets attached.
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "08.09.2023"
$ws.Range("I11").Style = "Normal"

$ws.Range("J11").Value = "17:38:06"

# K11 is an empty text cell (t="str" / shared-string "").  A plain "" clears
# the cell entirely, so use a text-prefix apostrophe to force an empty
# string value, then reset the style to drop the quote-prefix styling.
$ws.Range("K11").Value = "'"
$ws.Range("K11").Style = "Normal"

# --- Row 12 ---
$ws.Range("A12").Value = 19
$ws.Range("B12").Value = 162
$ws.Range("C12").Value = "09.09.2023, 22:59:56"
$ws.Range("D12").Value = "Старый город"
$ws.Range("E12").Value = "оливковый"
$ws.Range("F12").Value = "УЦЕНКА"
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = "Продукция Плиточка/Тротуарная плитка/Вибропресс/Старый город"

$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "09.09.2023"
$ws.Range("I12").Style = "Normal"

$ws.Range("J12").Value = "22:59:42"

$ws.Range("K12").Value = "'"
$ws.Range("K12").Style = "Normal"
